# Update betting-odds figures on Sheet1 for the Jogos da Semana workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("Q2").Value = 1.9
$ws.Range("R2").Value = 1.95
$ws.Range("S2").Value = 3.2
$ws.Range("T2").Value = 1.36

# Row 3
$ws.Range("H3").Value = 3.9
$ws.Range("I3").Value = 1.53
$ws.Range("K3").Value = 2.3
$ws.Range("Q3").Value = 1.88
$ws.Range("R3").Value = 1.98
$ws.Range("S3").Value = 3.25
$ws.Range("T3").Value = 1.33
$ws.Range("W3").Value = 1.83
$ws.Range("X3").Value = 1.83
$ws.Range("Y3").Value = 17
$ws.Range("Z3").Value = 34
$ws.Range("AD3").Value = 41
$ws.Range("AE3").Value = 11
$ws.Range("AI3").Value = 7

# Row 4
$ws.Range("G4").Value = 1.53
$ws.Range("W4").Value = 1.8
$ws.Range("X4").Value = 1.91
$ws.Range("Y4").Value = 7.5
$ws.Range("Z4").Value = 7.5
$ws.Range("AB4").Value = 11
$ws.Range("AG4").Value = 17
$ws.Range("AO4").Value = 251

# Row 5
$ws.Range("K5").Value = 3.65
$ws.Range("Q5").Value = 1.17
$ws.Range("R5").Value = 4.6
$ws.Range("Z5").Value = 7.3
$ws.Range("AA5").Value = 13.5
$ws.Range("AE5").Value = 30
$ws.Range("AN5").Value = 200

# Row 7
$ws.Range("N7").Value = 26
$ws.Range("O7").Value = 1.07
$ws.Range("P7").Value = 9
$ws.Range("Q7").Value = 1.25
$ws.Range("R7").Value = 4
$ws.Range("S7").Value = 1.67
$ws.Range("T7").Value = 2.2
